$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1419.4
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 6725.3335
$ws.Range("I31").Value = 6725.3335
$ws.Range("K31").Value = 20176.0005
$ws.Range("M31").Value = -19946.0005
$ws.Range("H86").Value = 2199.2727
$ws.Range("I86").Value = 1877.8
$ws.Range("K86").Value = 1877.8
$ws.Range("M86").Value = -754.8
$ws.Range("H88").Value = 1739.4445
$ws.Range("I88").Value = 789.3333
$ws.Range("J88").Value = 2214.5
$ws.Range("K88").Value = 789.3333
$ws.Range("L88").Value = 2214.5
$ws.Range("M88").Value = -383.3333
$ws.Range("N88").Value = -3026.5
$ws.Range("H89").Value = 2199.2727
$ws.Range("I89").Value = 1877.8
$ws.Range("K89").Value = 9389
$ws.Range("M89").Value = -3773
$ws.Range("H91").Value = 1739.4445
$ws.Range("I91").Value = 789.3333
$ws.Range("J91").Value = 2214.5
$ws.Range("K91").Value = 789.3333
$ws.Range("L91").Value = 2214.5
$ws.Range("M91").Value = 614.6667
$ws.Range("N91").Value = -5022.5
$ws.Range("H107").Value = 2502.25
$ws.Range("I107").Value = 2502.25
$ws.Range("K107").Value = 2502.25
$ws.Range("M107").Value = -582.25
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4876.923
$ws.Range("I32").Value = 3422.3777
$ws.Range("K32").Value = 3422.3777
$ws.Range("M32").Value = -3135.3777
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H74").Value = 625.65
$ws.Range("I74").Value = 625.65
$ws.Range("K74").Value = 625.65
$ws.Range("M74").Value = 248.35
$ws.Range("H77").Value = 625.65
$ws.Range("I77").Value = 625.65
$ws.Range("K77").Value = 3128.25
$ws.Range("M77").Value = 1239.75
$ws.Range("H88").Value = 1969.75
$ws.Range("I88").Value = 900
$ws.Range("J88").Value = 2326.3333
$ws.Range("K88").Value = 900
$ws.Range("L88").Value = 2326.3333
$ws.Range("M88").Value = -494
$ws.Range("N88").Value = -3138.3333
$ws.Range("H91").Value = 1969.75
$ws.Range("I91").Value = 900
$ws.Range("J91").Value = 2326.3333
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 2326.3333
$ws.Range("M91").Value = 504
$ws.Range("N91").Value = -5134.3333

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5500
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877
$ws.Range("H89").Value = 5500
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384
$ws.Range("H105").Value = 3144.6667
$ws.Range("I105").Value = 2717.25
$ws.Range("K105").Value = 2717.25
$ws.Range("M105").Value = -970.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H69").Value = 21996.25
$ws.Range("I69").Value = 20995
$ws.Range("J69").Value = 25000
$ws.Range("K69").Value = 20995
$ws.Range("L69").Value = 25000
$ws.Range("M69").Value = -20246
$ws.Range("N69").Value = -26498
$ws.Range("H72").Value = 21996.25
$ws.Range("I72").Value = 20995
$ws.Range("J72").Value = 25000
$ws.Range("K72").Value = 62985
$ws.Range("L72").Value = 75000
$ws.Range("M72").Value = -59241
$ws.Range("N72").Value = -82488
$ws.Range("H94").Value = 558.7778
$ws.Range("I94").Value = 550
$ws.Range("K94").Value = 550
$ws.Range("M94").Value = -99
$ws.Range("H105").Value = 6030.45
$ws.Range("I105").Value = 1222.2
$ws.Range("K105").Value = 1222.2
$ws.Range("M105").Value = 524.8
$ws.Range("H132").Value = 1577.9474
$ws.Range("I132").Value = 1379.8462
$ws.Range("J132").Value = 2007.1666
$ws.Range("K132").Value = 4139.5386
$ws.Range("L132").Value = 6021.4998
$ws.Range("M132").Value = -1609.5386
$ws.Range("N132").Value = -11081.4998
$ws.Range("H134").Value = 2406.1082
$ws.Range("I134").Value = 2216.4285
$ws.Range("K134").Value = 6649.2855
$ws.Range("M134").Value = -4114.2855

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3279.9
$ws.Range("I55").Value = 2199.75
$ws.Range("J55").Value = 4000
$ws.Range("K55").Value = 6599.25
$ws.Range("L55").Value = 12000
$ws.Range("M55").Value = -6422.25
$ws.Range("N55").Value = -12354
$ws.Range("H104").Value = 7500
$ws.Range("I104").Value = 5000
$ws.Range("J104").Value = 10000
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 30000
$ws.Range("M104").Value = -12379
$ws.Range("N104").Value = -35242
$ws.Range("H131").Value = 990
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 990
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2970
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -13050

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1194.2222
$ws.Range("I97").Value = 1256
$ws.Range("K97").Value = 1256
$ws.Range("M97").Value = -760

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5860.125
$ws.Range("J22").Value = 7536.8
$ws.Range("L22").Value = 7536.8
$ws.Range("N22").Value = -8126.8
$ws.Range("H27").Value = 5860.125
$ws.Range("J27").Value = 7536.8
$ws.Range("L27").Value = 7536.8
$ws.Range("N27").Value = -7750.8
$ws.Range("H50").Value = 12386.667

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 10125
$ws.Range("I54").Value = 6571.4287
$ws.Range("J54").Value = 35000
$ws.Range("K54").Value = 6571.4287
$ws.Range("L54").Value = 35000
$ws.Range("M54").Value = -6051.4287
$ws.Range("N54").Value = -36040
$ws.Range("H62").Value = 7088.722
$ws.Range("I62").Value = 5599.25
$ws.Range("J62").Value = 7514.2856
$ws.Range("K62").Value = 5599.25
$ws.Range("L62").Value = 7514.2856
$ws.Range("M62").Value = -4975.25
$ws.Range("N62").Value = -8762.285599999999
$ws.Range("H65").Value = 7088.722
$ws.Range("I65").Value = 5599.25
$ws.Range("J65").Value = 7514.2856
$ws.Range("K65").Value = 27996.25
$ws.Range("L65").Value = 37571.428
$ws.Range("M65").Value = -24876.25
$ws.Range("N65").Value = -43811.428
$ws.Range("H122").Value = 2845
$ws.Range("I122").Value = 2845
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8535
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6085
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1808.7858
$ws.Range("I126").Value = 1866.7273
$ws.Range("K126").Value = 5600.1819
$ws.Range("M126").Value = -3130.1819
